# Update values in "Random Forest Results Test Data" and
# "Random Forest Results All Data" sheets, and delete the last three
# rows (PC1_Dist_Prod, PC2_Dist_Prod, Avg_Distancia_casas_parcelas) from
# each, per the commit "Update-Delete of Sumfone Tagadad".

$wb = $excel.ActiveWorkbook

$testData = $wb.Worksheets.Item("Random Forest Results Test Data")
$allData  = $wb.Worksheets.Item("Random Forest Results All Data")

# New B/C/D values for rows 2-12 on the "Test Data" sheet.
$testValues = @(
    @(0.5980079174311927, 3.600506908756043, 2.812692307692308),
    @(-0.2629597305389222, 1.579881933271257, 1.240384615384615),
    @(0.09442647975077867, 1.311508233353434, 1.077692307692308),
    @(-0.3188192913385828, 0.995522669061526, 0.8203846153846155),
    @(-0.5728055401662053, 0.9164689763518373, 0.6515384615384616),
    @(-0.1774050408719345, 1.599015322002888, 1.165),
    @(-0.001158297413793141, 1.65793314331197, 1.338846153846154),
    @(0.1013091975308642, 2.075411914629146, 1.557307692307692),
    @(-0.1686827385892113, 2.886678954198991, 2.339230769230769),
    @(0.09974291938997792, 3.12734834841079, 2.450769230769231),
    @(-0.1812282258064515, 0.6582932477247506, 0.4988461538461539)
)

# New B/C/D values for rows 2-12 on the "All Data" sheet.
$allValues = @(
    @(0.8164400018691669, 2.278854601455223, 1.613643410852713),
    @(0.5022008025682183, 0.7960142961592602, 0.5075968992248062),
    @(0.7312563508808947, 0.7825113615838124, 0.5840310077519381),
    @(0.6552904873752202, 0.5939422883247708, 0.4131782945736434),
    @(0.7559024387424873, 0.6170365536393543, 0.4290697674418604),
    @(0.7429353125114289, 0.919114420422973, 0.6383720930232558),
    @(0.6856788805577479, 0.9281894520596152, 0.6481395348837209),
    @(0.7724230241296854, 1.327208412625024, 0.852248062015504),
    @(0.7009614902198025, 1.63505185656095, 1.161317829457364),
    @(0.7285962674440363, 1.819908848556505, 1.32108527131783),
    @(0.750996020466174, 0.3973946157412743, 0.2685271317829457)
)

for ($i = 0; $i -lt $testValues.Count; $i++) {
    $row = $i + 2
    $testData.Cells.Item($row, 2).Value = $testValues[$i][0]
    $testData.Cells.Item($row, 3).Value = $testValues[$i][1]
    $testData.Cells.Item($row, 4).Value = $testValues[$i][2]
}

for ($i = 0; $i -lt $allValues.Count; $i++) {
    $row = $i + 2
    $allData.Cells.Item($row, 2).Value = $allValues[$i][0]
    $allData.Cells.Item($row, 3).Value = $allValues[$i][1]
    $allData.Cells.Item($row, 4).Value = $allValues[$i][2]
}

# Delete rows 13-15 (PC1_Dist_Prod, PC2_Dist_Prod, Avg_Distancia_casas_parcelas)
# on both sheets.
$testData.Range("A13:D15").Delete()
$allData.Range("A13:D15").Delete()
